$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REFERENCES")

$ws.Range("J2").Value = "PCS per PU"
$ws.Range("K2").Value = "PU per HU"

[void]$ws.Range("K3").Select()
